# Remove the "select_one / country_1 (autocomplete datalist)" demo row
# from the survey sheet (fix change detection for custom date pickers:
# the datalist-based select row is unrelated to the custom date picker
# prompts and gets dropped as part of this cleanup). All subsequent rows
# shift up by one.

$wb = $excel.ActiveWorkbook

$wsSurvey = $wb.Worksheets.Item("survey")
$wsPromptTypes = $wb.Worksheets.Item("prompt_types")

# Update the selection on the "prompt_types" sheet first (select whole
# column C), then return focus to "survey" so it stays the active tab.
$wsPromptTypes.Activate()
$wsPromptTypes.Range("C1:C1048576").Select()

$wsSurvey.Activate()

# Delete row 5 (type=select_one, name=country_1, the "Select from list
# (autocomplete datalist)" prompt) - everything below shifts up.
$wsSurvey.Rows("5:5").Delete()

# Reset the view: no more frozen/scrolled topLeftCell, and select C4
# (now the "if" clause row that used to be row 6).
$wsSurvey.Range("C4").Select()
